$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 3.537761648806719
$ws.Range("E2").Value = 0.4942365360607697
$ws.Range("G2").Value = 8.974608811992548

$ws.Range("B3").Value = 0.6606524410359556
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 6.348428708163715

$ws.Range("B4").Value = 1.455362044514542
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.7527432677738641
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 14.05633640148523

$ws.Range("B5").Value = 1.455362044514542
$ws.Range("C5").Value = 3286.919754855326
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 10.19245300693656
$ws.Range("G5").Value = 3299.320313174551

$ws.Range("B6").Value = 3.286832544864788
$ws.Range("C6").Value = 1.655778082260271
$ws.Range("D6").Value = 0.7527432677738641
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("G6").Value = 6.189590430959694

$ws.Range("B7").Value = 3.286832544864788
$ws.Range("C7").Value = 1.655778082260271
$ws.Range("D7").Value = 0.7527432677738641
$ws.Range("E7").Value = 0.4942365360607697
$ws.Range("G7").Value = 6.189590430959694

$ws.Range("B8").Value = 3.286832544864788
$ws.Range("C8").Value = 1.655778082260271
$ws.Range("D8").Value = 22.3905356188092
$ws.Range("E8").Value = 0.4942365360607697
$ws.Range("G8").Value = 27.82738278199502

$ws.Range("B9").Value = 0.1190320826869504
$ws.Range("C9").Value = 0.306821227259698
$ws.Range("D9").Value = 0.1494219747398047
$ws.Range("E9").Value = 0.4942365360607697
$ws.Range("G9").Value = 1.069511820747223

$ws.Range("B10").Value = 3.286832544864788
$ws.Range("C10").Value = 1.655778082260271
$ws.Range("D10").Value = 3.537761648806719
$ws.Range("E10").Value = 0.4942365360607697
$ws.Range("G10").Value = 8.974608811992548

$ws.Range("B11").Value = 0.6606524410359556
$ws.Range("C11").Value = 0.306821227259698
$ws.Range("D11").Value = 0.7527432677738641
$ws.Range("E11").Value = 10.19245300693656
$ws.Range("G11").Value = 11.91266994300607

$ws.Range("B12").Value = 0.2917716402565462
$ws.Range("C12").Value = 10.34677158129881
$ws.Range("D12").Value = 0.7527432677738641
$ws.Range("E12").Value = 10.19245300693656
$ws.Range("G12").Value = 21.58373949626578

$ws.Range("B13").Value = 1.455362044514542
$ws.Range("C13").Value = 0.306821227259698
$ws.Range("D13").Value = 3.537761648806719
$ws.Range("E13").Value = 0.4942365360607697
$ws.Range("G13").Value = 5.794181456641729

$ws.Range("B14").Value = 0.04271373187048222
$ws.Range("C14").Value = 0.306821227259698
$ws.Range("D14").Value = 0.1494219747398047
$ws.Range("E14").Value = 0.4942365360607697
$ws.Range("G14").Value = 0.9931934699307545

$ws.Range("B15").Value = 1.455362044514542
$ws.Range("C15").Value = 1.655778082260271
$ws.Range("D15").Value = 3.537761648806719
$ws.Range("E15").Value = 0.4942365360607697
$ws.Range("G15").Value = 7.143138311642302
